# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "68.650.82"
$ws.Range("E2").Value = "  -2.51%  "
$ws.Range("D3").Value = "3.453.76"
$ws.Range("E3").Value = "  -4.64%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "574.53"
$ws.Range("E5").Value = "  -4.56%  "
$ws.Range("E6").Value = "  -3.25%  "
$ws.Range("D7").Value = "3.444.18"
$ws.Range("E7").Value = "  -4.59%  "
$ws.Range("E8").Value = "  -3.78%  "
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("E10").Value = "  -4.64%  "
$ws.Range("D11").Value = "0.612"
$ws.Range("E11").Value = "  -5.43%  "
$ws.Range("D12").Value = "50.86"
$ws.Range("E12").Value = "  -4.47%  "
$ws.Range("E13").Value = "  -6.65%  "
$ws.Range("D14").Value = "9.03"
$ws.Range("E14").Value = "  -5.76%  "
$ws.Range("D15").Value = "3.998.76"
$ws.Range("E15").Value = "  -4.75%  "
$ws.Range("D16").Value = "635.04"
$ws.Range("E16").Value = "  +4.65%  "
$ws.Range("D17").Value = "68.505.05"
$ws.Range("E17").Value = "  -2.80%  "
$ws.Range("D18").Value = "3.456.04"
$ws.Range("E18").Value = "  -4.56%  "
$ws.Range("D19").Value = "12.27"
$ws.Range("E19").Value = "  -4.96%  "
$ws.Range("E20").Value = "  -2.47%  "
$ws.Range("D21").Value = "18.08"
$ws.Range("E21").Value = "  -5.43%  "
$ws.Range("D22").Value = "0.936"
$ws.Range("E22").Value = "  -6.38%  "
$ws.Range("D23").Value = "17.83"
$ws.Range("E23").Value = "  -2.11%  "
$ws.Range("D24").Value = "5.32"
$ws.Range("E24").Value = "  +2.59%  "
$ws.Range("D25").Value = "99.06"
$ws.Range("E25").Value = "  -3.92%  "
$ws.Range("D26").Value = "4.26"
$ws.Range("E26").Value = "  -7.82%  "
$ws.Range("E27").Value = "  -6.05%  "
$ws.Range("E28").Value = "  +1.75%  "
$ws.Range("D29").Value = "9.74"
$ws.Range("E29").Value = "  -8.65%  "
$ws.Range("D30").Value = "9.17"
$ws.Range("E30").Value = "  -5.28%  "
$ws.Range("D31").Value = "32.23"
$ws.Range("E31").Value = "  -4.70%  "
$ws.Range("D32").Value = "4.16"
$ws.Range("E32").Value = "  -11.46%  "
$ws.Range("D33").Value = "6.69"
$ws.Range("E33").Value = "  -8.57%  "
$ws.Range("D34").Value = "11.53"
$ws.Range("E34").Value = "  -6.12%  "
$ws.Range("D35").Value = "60.96"
$ws.Range("E35").Value = "  -3.99%  "
$ws.Range("E36").Value = "  -7.78%  "
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("D38").Value = "3.657.13"
$ws.Range("E38").Value = "  -6.33%  "
$ws.Range("D39").Value = "500.16"
$ws.Range("E39").Value = "  -3.31%  "
$ws.Range("D40").Value = "0.0₃0772"
$ws.Range("E40").Value = "  -12.41%  "
$ws.Range("D41").Value = "3.50"
$ws.Range("E41").Value = "  -1.12%  "
$ws.Range("D42").Value = "2.86"
$ws.Range("E42").Value = "  -6.67%  "
$ws.Range("D43").Value = "0.365"
$ws.Range("E43").Value = "  -5.94%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").Value = "34.21"
$ws.Range("E44").Value = "  -7.33%  "
$ws.Range("B45").Value = "Kaspa"
$ws.Range("C45").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D45").Value = "0.132"
$ws.Range("E45").Value = "  -2.20%  "
$ws.Range("D46").Value = "3.40"
$ws.Range("E46").Value = "  +63.20%  "
$ws.Range("D47").Value = "0.0434"
$ws.Range("E47").Value = "  -5.64%  "
$ws.Range("E48").Value = "  -5.84%  "
$ws.Range("E49").Value = "  -3.96%  "
$ws.Range("E50").Value = "  -4.82%  "
$ws.Range("D51").Value = "0.999"
$ws.Range("E51").Value = "  -0.36%  "
